$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "94.663.20"
Set-TextValue "E2" "  -2.17%  "

Set-TextValue "D3" "3.467.00"
Set-TextValue "E3" "  +3.88%  "

Set-TextValue "E4" "  +0.09%  "

Set-TextValue "D5" "237.82"
Set-TextValue "E5" "  -4.94%  "

Set-TextValue "D6" "638.97"
Set-TextValue "E6" "  -2.48%  "

Set-TextValue "E7" "  +1.94%  "

Set-TextValue "D8" "0.399"
Set-TextValue "E8" "  -5.53%  "

Set-TextValue "D10" "0.979"
Set-TextValue "E10" "  -2.78%  "

Set-TextValue "D11" "3.462.85"
Set-TextValue "E11" "  +3.84%  "

Set-TextValue "D12" "42.22"
Set-TextValue "E12" "  +3.85%  "

Set-TextValue "D13" "0.197"
Set-TextValue "E13" "  -4.40%  "

Set-TextValue "D14" "6.19"
Set-TextValue "E14" "  +1.65%  "

Set-TextValue "B15" "WrappedliquidstakedEther2.0"
Set-TextValue "C15" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D15" "4.128.71"
Set-TextValue "E15" "  +4.23%  "

Set-TextValue "B16" "WrappedBTC"
Set-TextValue "C16" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D16" "94.626.31"
Set-TextValue "E16" "  -1.95%  "

Set-TextValue "D17" "0.0000254"
Set-TextValue "E17" "  +0.75%  "

Set-TextValue "D18" "8.40"
Set-TextValue "E18" "  -3.33%  "

Set-TextValue "D19" "3.487.28"
Set-TextValue "E19" "  +4.37%  "

Set-TextValue "D20" "17.71"
Set-TextValue "E20" "  +1.58%  "

Set-TextValue "D21" "11.40"
Set-TextValue "E21" "  +6.86%  "

Set-TextValue "D22" "0.501"
Set-TextValue "E22" "  -12.22%  "

Set-TextValue "D23" "500.82"
Set-TextValue "E23" "  -1.54%  "

Set-TextValue "D24" "3.16"
Set-TextValue "E24" "  -5.23%  "

Set-TextValue "D25" "6.67"
Set-TextValue "E25" "  +0.79%  "

Set-TextValue "D26" "0.0000191"
Set-TextValue "E26" "  -3.93%  "

Set-TextValue "D27" "91.53"
Set-TextValue "E27" "  -5.15%  "

Set-TextValue "D28" "3.652.67"
Set-TextValue "E28" "  +3.82%  "

Set-TextValue "D29" "12.00"
Set-TextValue "E29" "  -0.86%  "

Set-TextValue "D30" "11.79"
Set-TextValue "E30" "  +3.60%  "

Set-TextValue "E31" "  +0.17%  "

Set-TextValue "D32" "2.74"
Set-TextValue "E32" "  +8.36%  "

Set-TextValue "E33" "  -6.72%  "

Set-TextValue "E34" "  -2.96%  "

Set-TextValue "D35" "0.998"
Set-TextValue "E35" "  -0.01%  "

Set-TextValue "D36" "30.40"
Set-TextValue "E36" "  +7.22%  "

Set-TextValue "D37" "0.568"
Set-TextValue "E37" "  +2.71%  "

Set-TextValue "D38" "540.20"
Set-TextValue "E38" "  +6.74%  "

Set-TextValue "D39" "7.68"
Set-TextValue "E39" "  -2.05%  "

Set-TextValue "D40" "1.45"
Set-TextValue "E40" "  -3.91%  "

Set-TextValue "D41" "0.937"
Set-TextValue "E41" "  +12.12%  "

Set-TextValue "E42" "  +0.60%  "

Set-TextValue "E43" "  -0.01%  "

Set-TextValue "D44" "24.05"
Set-TextValue "E44" "  -1.25%  "

Set-TextValue "E45" "  +1.11%  "

Set-TextValue "D46" "5.58"
Set-TextValue "E46" "  +0.29%  "

Set-TextValue "D47" "0.0412"
Set-TextValue "E47" "  -4.84%  "

Set-TextValue "D48" "2.16"
Set-TextValue "E48" "  +9.10%  "

Set-TextValue "E49" "  -5.07%  "

Set-TextValue "B50" "OKB"
Set-TextValue "C50" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D50" "53.51"
Set-TextValue "E50" "  -2.09%  "

Set-TextValue "B51" "dogwifhat"
Set-TextValue "C51" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D51" "3.20"
Set-TextValue "E51" "  +3.13%  "
